$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.084007382392883
$ws.Range("B1").Value = 2.414854049682617
$ws.Range("C1").Value = 6.449297428131104
$ws.Range("D1").Value = 2.220097541809082
$ws.Range("E1").Value = 1.2778400182724
